$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append below the existing data (rows 227-229),
# continuing the daily series in column A (dates) / B / C / D.
$data = @(
    @{ Row = 227; A = 44301; B = 2; C = 24; D = 150.281778334377 },
    @{ Row = 228; A = 44302; B = 0; C = 19; D = 118.9730745147151 },
    @{ Row = 229; A = 44303; B = 2; C = 17; D = 106.4495929868503 }
)

foreach ($item in $data) {
    $r = $item.Row

    # Column A uses a special style (bold, centered, bordered, custom date/time
    # number format) applied to every date cell in the column. Copy the format
    # from the row directly above so the new cell reuses the existing style
    # entry instead of creating a new one.
    $ws.Range("A$($r - 1)").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}

$excel.CutCopyMode = 0
